# Update Cadm1-Crtam LR-pairs sheet with new TPM-based values.
# The within-cluster (self-pair) rows are removed and the remaining
# sending/target cluster combinations get recalculated TPM-derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('ECs', 'Cadm1', 'Crtam', 'FAPs', '3', '1', '0.8244693333333334', '2.473408', '0.04531162661728313', '0.04531162661728313', '2', '0.6666666666666666', '0.4937903333333333', '1.481371', '0.8374020920199343', '0.8374020920199343', '0.4071149869297778', '3.664034882368', '0.03794405092213903', '0.03794405092213903'),
    @('ECs', 'Cadm1', 'Crtam', 'MuSCs', '3', '1', '0.8244693333333334', '2.473408', '0.04531162661728313', '0.04531162661728313', '2', '0.6666666666666666', '0.09091933333333334', '0.272758', '0.1541869793692284', '0.1541869793692284', '0.07496020214044445', '0.674641819264', '0.006986462838425214', '0.006986462838425214'),
    @('ECs', 'Cadm1', 'Crtam', 'Resolving-Mac', '3', '1', '0.8244693333333334', '2.473408', '0.04531162661728313', '0.04531162661728313', '1', '0.3333333333333333', '0.004959666666666667', '0.014879', '0.008410928610837261', '0.008410928610837259', '0.004089093070222222', '0.036801837632', '0.0003811128567188818', '0.0003811128567188817'),
    @('FAPs', 'Cadm1', 'Crtam', 'FAPs', '1', '0.3333333333333333', '0.01719833333333333', '0.051595', '0.0009451952024569835', '0.0009451952024569837', '2', '0.6666666666666666', '0.4937903333333333', '1.481371', '0.8374020920199343', '0.8374020920199343', '0.008492370749444445', '0.07643133674500001', '0.0007915084399046833', '0.0007915084399046835'),
    @('FAPs', 'Cadm1', 'Crtam', 'MuSCs', '1', '0.3333333333333333', '0.01719833333333333', '0.051595', '0.0009451952024569835', '0.0009451952024569837', '2', '0.6666666666666666', '0.09091933333333334', '0.272758', '0.1541869793692284', '0.1541869793692284', '0.001563661001111111', '0.01407294901', '0.0001457367931811286', '0.0001457367931811286'),
    @('FAPs', 'Cadm1', 'Crtam', 'Resolving-Mac', '1', '0.3333333333333333', '0.01719833333333333', '0.051595', '0.0009451952024569835', '0.0009451952024569837', '1', '0.3333333333333333', '0.004959666666666667', '0.014879', '0.008410928610837261', '0.008410928610837259', '8.529800055555556E-05', '0.000767682005', '7.949969371171559E-06', '7.949969371171559E-06'),
    @('MuSCs', 'Cadm1', 'Crtam', 'FAPs', '3', '1', '3.182317', '9.546951', '0.1748954798583565', '0.1748954798583565', '2', '0.6666666666666666', '0.4937903333333333', '1.481371', '0.8374020920199343', '0.8374020920199343', '1.571397372202333', '14.142576349821', '0.146457840718218', '0.146457840718218'),
    @('MuSCs', 'Cadm1', 'Crtam', 'MuSCs', '3', '1', '3.182317', '9.546951', '0.1748954798583565', '0.1748954798583565', '2', '0.6666666666666666', '0.09091933333333334', '0.272758', '0.1541869793692284', '0.1541869793692284', '0.2893341400953333', '2.604007260858', '0.02696660574469171', '0.02696660574469171'),
    @('MuSCs', 'Cadm1', 'Crtam', 'Resolving-Mac', '3', '1', '3.182317', '9.546951', '0.1748954798583565', '0.1748954798583565', '1', '0.3333333333333333', '0.004959666666666667', '0.014879', '0.008410928610837261', '0.008410928610837259', '0.01578323154766667', '0.142049083929', '0.001471033395446762', '0.001471033395446762'),
    @('Resolving-Mac', 'Cadm1', 'Crtam', 'FAPs', '3', '1', '14.17155133333333', '42.514654', '0.7788476983219034', '0.7788476983219035', '2', '0.6666666666666666', '0.4937903333333333', '1.481371', '0.8374020920199343', '0.8374020920199343', '6.997775056737111', '62.979975510634', '0.6522086919396726', '0.6522086919396727'),
    @('Resolving-Mac', 'Cadm1', 'Crtam', 'MuSCs', '3', '1', '14.17155133333333', '42.514654', '0.7788476983219034', '0.7788476983219035', '2', '0.6666666666666666', '0.09091933333333334', '0.272758', '0.1541869793692284', '0.1541869793692284', '1.288467999525778', '11.596211995732', '0.1200881739929303', '0.1200881739929304'),
    @('Resolving-Mac', 'Cadm1', 'Crtam', 'Resolving-Mac', '3', '1', '14.17155133333333', '42.514654', '0.7788476983219034', '0.7788476983219035', '1', '0.3333333333333333', '0.004959666666666667', '0.014879', '0.008410928610837261', '0.008410928610837259', '0.07028617076288889', '0.632575536866', '0.006550832389300445', '0.006550832389300445'),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    # Columns A-D are text (cluster/gene names)
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    # Columns E-T are numeric
    for ($c = 5; $c -le 20; $c++) {
        $ws.Cells.Item($row, $c).Value = [double]$vals[$c - 1]
    }
}

# Original sheet had 16 data rows (rows 2-17); the new TPM data only has
# 12 data rows (rows 2-13). Remove the now-obsolete trailing rows.
$ws.Rows("14:17").Delete()
